# Sweden Allsvenskan update (2024-01-29 18:52) — corrects mismatched/duplicated
# match rows: several consecutive rows had their HomeTeam/AwayTeam/odds data
# swapped between the wrong "id" rows. This script restores the correct
# pairing by swapping (or, for one group, cyclically rotating) the B:AC
# payload between the affected rows while leaving column A (the row's
# sequential id) untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowData($row) {
    $b = $ws.Range("B$row").Value()
    $fac = $ws.Range("F$row`:AC$row").Value()
    return @{ B = $b; FAC = $fac }
}

function Set-RowData($row, $data) {
    $ws.Range("B$row").Value = $data.B
    $ws.Range("F$row`:AC$row").Value = $data.FAC
}

# ---- simple pairwise swaps (B and F:AC payload exchanged between the two rows) ----
$swapPairs = @(
    @(304, 305),
    @(350, 352),
    @(450, 452),
    @(518, 519),
    @(523, 525),
    @(597, 598),
    @(609, 610)
)

foreach ($pair in $swapPairs) {
    $r1 = $pair[0]
    $r2 = $pair[1]
    $d1 = Get-RowData $r1
    $d2 = Get-RowData $r2
    Set-RowData $r1 $d2
    Set-RowData $r2 $d1
}

# ---- cyclic rotation: rows 461-464 each take on the payload that used to
# belong to the previous row in the list (461 wraps around and takes 464's) ----
$cycleRows = @(461, 462, 463, 464)
$cycleData = @()
foreach ($r in $cycleRows) {
    $cycleData += ,(Get-RowData $r)
}

$n = $cycleRows.Count
for ($i = 0; $i -lt $n; $i++) {
    $srcIdx = ($i - 1 + $n) % $n
    Set-RowData $cycleRows[$i] $cycleData[$srcIdx]
}
